$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 413
$ws.Range("F3").Value = 1089
$ws.Range("F4").Value = 9695
$ws.Range("F5").Value = 210
$ws.Range("F8").Value = 6604
$ws.Range("F10").Value = 10515
$ws.Range("F11").Value = 11687
$ws.Range("F12").Value = 1261
$ws.Range("F13").Value = 1213
$ws.Range("F14").Value = 5068
$ws.Range("F15").Value = 841
$ws.Range("F16").Value = 503
$ws.Range("F18").Value = 344
$ws.Range("F19").Value = 184
$ws.Range("F21").Value = 285
$ws.Range("F22").Value = 1920
$ws.Range("F23").Value = 929
$ws.Range("F24").Value = 1332
$ws.Range("F26").Value = 8
$ws.Range("F27").Value = 2082
$ws.Range("F28").Value = 446
$ws.Range("F29").Value = 673
$ws.Range("F30").Value = 2765
$ws.Range("F31").Value = 207
$ws.Range("F32").Value = 1848
$ws.Range("F33").Value = 98
$ws.Range("F34").Value = 840
$ws.Range("F35").Value = 92
$ws.Range("F36").Value = 942
$ws.Range("F37").Value = 38
$ws.Range("F38").Value = 58
$ws.Range("F39").Value = 3450
$ws.Range("F41").Value = 99
$ws.Range("F42").Value = 534
$ws.Range("F43").Value = 604
$ws.Range("F45").Value = 906
$ws.Range("F47").Value = 12
$ws.Range("F48").Value = 4244
$ws.Range("F49").Value = 97

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 14
$ws.Range("F9").Value = 35
$ws.Range("F12").Value = 69
$ws.Range("F26").Value = 46
$ws.Range("F28").Value = 60

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 6123

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 413
$ws.Range("F3").Value = 1089
$ws.Range("F4").Value = 9695
$ws.Range("F6").Value = 14
$ws.Range("F9").Value = 10515
$ws.Range("F10").Value = 11687
$ws.Range("F12").Value = 1213
$ws.Range("F13").Value = 5068
$ws.Range("F14").Value = 841
$ws.Range("F15").Value = 503
$ws.Range("F17").Value = 344
$ws.Range("F18").Value = 35
$ws.Range("F19").Value = 184
$ws.Range("F21").Value = 285
$ws.Range("F22").Value = 1920
$ws.Range("F23").Value = 929
$ws.Range("F24").Value = 1332
$ws.Range("F26").Value = 2082
$ws.Range("F27").Value = 446
$ws.Range("F28").Value = 673
$ws.Range("F29").Value = 2765
$ws.Range("F30").Value = 207
$ws.Range("F31").Value = 1848
$ws.Range("F32").Value = 98
$ws.Range("F34").Value = 840
$ws.Range("F38").Value = 92
$ws.Range("F39").Value = 942
$ws.Range("F40").Value = 38
$ws.Range("F41").Value = 46
$ws.Range("F43").Value = 99
$ws.Range("F44").Value = 534
$ws.Range("F45").Value = 604
$ws.Range("F46").Value = 906
$ws.Range("F48").Value = 4244
$ws.Range("F49").Value = 97
